# "i update my skill" - mark skills 9..36 with a checkmark in the
# C++Primer sheet, and leave the workbook focused on that sheet/cell,
# mirroring the selection/view state recorded in the authored diff.

$wb = $excel.ActiveWorkbook

$wsPython = $wb.Worksheets.Item("python核心编程")
$wsCpp    = $wb.Worksheets.Item("C++Primer")

# Mark rows 9 through 36 (inclusive) in column C with a check mark,
# matching the newly added <c r="C9".."C36" t="s"><v>369</v></c> cells.
$wsCpp.Range("C9:C36").Value = "√"

# --- View / selection state -------------------------------------------------

# Sheet1 (python核心编程) is no longer the active tab; its remembered
# scroll position/selection moves to A99 / B2.
$wsPython.Activate()
$excel.ActiveWindow.ScrollRow = 99
$excel.ActiveWindow.ScrollColumn = 1
$wsPython.Range("B2").Select()

# Sheet3 (C++Primer) becomes the active tab, selection moves to G37.
$wsCpp.Activate()
$wsCpp.Range("G37").Select()
